$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Remember the width of the "type in mysql" column so the newly inserted column can match it.
$typeColWidth = $ws.Range("C1").ColumnWidth

# Insert a whole new column before column D ("is null?"), shifting is null?/is key?/initial
# value/description/mapped-from-MassBank-File one column to the right (D->E, E->F, F->G, G->H, H->I).
# This also duplicates formatting from the left-adjacent column (C) for the new column D,
# matching the source workbook's layout.
$ws.Columns("D:D").Insert()

# Give the new "type in JavaDB" column roughly the same width as "type in mysql".
$ws.Columns("D:D").ColumnWidth = $typeColWidth

# Populate the new "type in JavaDB" column header and its data cells.
$ws.Range("D4").Value = "type in JavaDB"
$ws.Range("D5").Value = "varchar(8)"
$ws.Range("D6").Value = "varchar(255) "
$ws.Range("D7").Value = "int"
$ws.Range("D8").Value = "int"

# Re-affirm the header cells that were shifted right, so the table reads their current text.
$ws.Range("E4").Value = "is null?"
$ws.Range("F4").Value = "is key?"
$ws.Range("G4").Value = "initial value"
$ws.Range("H4").Value = "description"

# Update the mysql type for the ION row (previously varchar(30)) to tinyint(4), reflecting the
# new JavaDB mapping work.
$ws.Range("C7").Value = "tinyint(4)"

# Grow the table (ListObject) so it covers the newly inserted column.
$lo.Resize($ws.Range("B4:I8"))

# Set the last header cell after resizing, so the table picks up the right-most column name.
$ws.Range("I4").Value = "mapped from MassBank File"

# Update the selection to match the post-edit cursor location.
$ws.Range("D9").Select() | Out-Null
